# Add team-member names to the "Library Management System" project's
# "Team Roles:" list (roles 1-8), per the commit:
# "adding readme for github sequence pushing and roles"

$d = $word.ActiveDocument
$count = $d.Paragraphs.Count

# Anchor on the "6. Library Management System" heading, then the
# "Team Roles:" paragraph that follows it, so we land on the right
# list even though similar headings ("Team Roles:", "N. CRUD Developer",
# "N. UI Developer", "N. Tester", ...) repeat for other projects earlier
# in the document.
$headingIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Library Management System*") {
        $headingIdx = $i
        break
    }
}

$teamRolesIdx = -1
for ($i = $headingIdx; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Team Roles:*") {
        $teamRolesIdx = $i
        break
    }
}

function Append-ParagraphText($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # $r.End points just past the trailing paragraph mark; back up one
    # character so the insertion lands before the mark (i.e. at the
    # visible end of the paragraph's text), not inside the next paragraph.
    $insertPoint = $d.Range($r.End - 1, $r.End - 1)
    # InsertBefore (rather than InsertAfter) correctly inherits the
    # formatting of the preceding run even when that run is immediately
    # followed by a zero-width marker (e.g. a spell-check proofErr) right
    # before the paragraph mark, as happens for "6. UI Developer" below.
    $insertPoint.InsertBefore($text)
}

Append-ParagraphText ($teamRolesIdx + 1) " – Mohamed Emad"
Append-ParagraphText ($teamRolesIdx + 2) " – Youssef Amr"
Append-ParagraphText ($teamRolesIdx + 3) " - Alkady "
Append-ParagraphText ($teamRolesIdx + 4) " – Mohamed Seif"
Append-ParagraphText ($teamRolesIdx + 5) " – Youssef amr"
Append-ParagraphText ($teamRolesIdx + 6) " – Mohamed Emad"
Append-ParagraphText ($teamRolesIdx + 7) " - Amr"
Append-ParagraphText ($teamRolesIdx + 8) " – Omar Salama"
